# ---------------------------------------------------------------------------
# C1--C2-and-C3-PowerPoint.pptx edit
#
# 1) The table on the "PLENARY - COMPLETE THE MISSING GAPS" slide (slide 16)
#    is switched from the deck's custom table style to the built-in
#    "Light Style 2" gallery style ({E1A88005-5D59-4A64-BADF-D830EDD0A735}).
# 2) The presentation's theme palette is switched from the custom "Integral"
#    colour scheme to the stock Office colour scheme (fonts / format scheme
#    are already identical between the two, only the 12 theme colours
#    change).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 16 --------------------------------------
$slide = $p.Slides.Item(16)

$tableShape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $candidate = $slide.Shapes.Item($i)
    if ($candidate.HasTable) {
        $tableShape = $candidate
        break
    }
}

if ($tableShape -ne $null) {
    $tableShape.Table.ApplyStyle("{E1A88005-5D59-4A64-BADF-D830EDD0A735}")
}

# --- 2. Recolour the theme: Integral -> Office ------------------------------
# PowerPoint's .RGB is a plain Long stored as 0x00BBGGRR (like VBA's RGB()),
# so build it by hand since this host has no RGB() helper.
function ToOleColor($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme

# index : scheme slot : old (Integral) -> new (Office) RGB value
$themeColors.Item(1).RGB  = ToOleColor 0x00 0x00 0x00   # dk1      000000
$themeColors.Item(2).RGB  = ToOleColor 0xFF 0xFF 0xFF   # lt1      FFFFFF
$themeColors.Item(3).RGB  = ToOleColor 0x44 0x54 0x6A   # dk2      44546A
$themeColors.Item(4).RGB  = ToOleColor 0xE7 0xE6 0xE6   # lt2      E7E6E6
$themeColors.Item(5).RGB  = ToOleColor 0x5B 0x9B 0xD5   # accent1  5B9BD5
$themeColors.Item(6).RGB  = ToOleColor 0xED 0x7D 0x31   # accent2  ED7D31
$themeColors.Item(7).RGB  = ToOleColor 0xA5 0xA5 0xA5   # accent3  A5A5A5
$themeColors.Item(8).RGB  = ToOleColor 0xFF 0xC0 0x00   # accent4  FFC000
$themeColors.Item(9).RGB  = ToOleColor 0x44 0x72 0xC4   # accent5  4472C4
$themeColors.Item(10).RGB = ToOleColor 0x70 0xAD 0x47   # accent6  70AD47
$themeColors.Item(11).RGB = ToOleColor 0x05 0x63 0xC1   # hlink    0563C1
$themeColors.Item(12).RGB = ToOleColor 0x95 0x4F 0x72   # folHlink 954F72
